$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 113, pushing the existing
# rows 113-118 down to 115-120.
$ws.Rows("113:114").Insert()

# Row 113: new "Primera" record (week of 2021-11-09)
$ws.Cells.Item(113, 1).Value = 8
$ws.Cells.Item(113, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(113, 3).Value = "Coquimbo"
$ws.Cells.Item(113, 4).Value = 44509
$ws.Cells.Item(113, 5).Value = 4
$ws.Cells.Item(113, 6).Value = 100112021
$ws.Cells.Item(113, 7).Value = "Ají"
$ws.Cells.Item(113, 8).Value = "Inferno"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 440
$ws.Cells.Item(113, 11).Value = 24000
$ws.Cells.Item(113, 12).Value = 25000
$ws.Cells.Item(113, 13).Value = 24500
$ws.Cells.Item(113, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(113, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(113, 16).Value = 2042
$ws.Cells.Item(113, 17).Value = 12
$ws.Cells.Item(113, 18).Value = "Hortaliza"

# Row 114: new "Segunda" record (week of 2021-11-09)
$ws.Cells.Item(114, 1).Value = 8
$ws.Cells.Item(114, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(114, 3).Value = "Coquimbo"
$ws.Cells.Item(114, 4).Value = 44509
$ws.Cells.Item(114, 5).Value = 4
$ws.Cells.Item(114, 6).Value = 100112021
$ws.Cells.Item(114, 7).Value = "Ají"
$ws.Cells.Item(114, 8).Value = "Inferno"
$ws.Cells.Item(114, 9).Value = "Segunda"
$ws.Cells.Item(114, 10).Value = 300
$ws.Cells.Item(114, 11).Value = 14000
$ws.Cells.Item(114, 12).Value = 15000
$ws.Cells.Item(114, 13).Value = 14500
$ws.Cells.Item(114, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(114, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(114, 16).Value = 1208
$ws.Cells.Item(114, 17).Value = 12
$ws.Cells.Item(114, 18).Value = "Hortaliza"
